# Update the "想去人数" (want-to-go count) column for the sheets that list
# exhibition events: "展览" and "全部类型". Rows 2, 3 and 5 each get
# incremented by 1 (318->319, 237->238, 277->278).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 319
    $ws.Range("F3").Value = 238
    $ws.Range("F5").Value = 278
}
